# Add season-record columns (Wins / Losses / Ties) to the Pittsburgh 2014
# player table: three new header cells in row 1 (AD1:AF1) formatted like
# the existing header row, and the team's season record repeated down
# every data row (AD2:AF47).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the rest of the header row (bold, centered/top
# aligned, thin box border) by copying the format from an existing
# header cell.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data: the Pirates' 2014 season record (88-74-0) on every row ------
$lastRow = 47
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 88   # column AD = Wins
    $ws.Cells.Item($r, 31).Value = 74   # column AE = Losses
    $ws.Cells.Item($r, 32).Value = 0    # column AF = Ties
}
